$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '44.462.57'
Set-TextValue 2 5 '  +0.59%  '

Set-TextValue 3 4 '2.228.61'
Set-TextValue 3 5 '  -0.77%  '

Set-TextValue 4 5 '  +0.46%  '

Set-TextValue 5 4 '303.51'
Set-TextValue 5 5 '  -0.95%  '

Set-TextValue 6 4 '93.73'
Set-TextValue 6 5 '  -2.30%  '

Set-TextValue 7 4 '0.567'
Set-TextValue 7 5 '  -1.10%  '

Set-TextValue 8 4 '1.01'
Set-TextValue 8 5 '  +0.23%  '

Set-TextValue 9 4 '0.512'
Set-TextValue 9 5 '  -3.10%  '

Set-TextValue 10 4 '34.44'
Set-TextValue 10 5 '  -2.14%  '

Set-TextValue 11 4 '0.0797'
Set-TextValue 11 5 '  -2.27%  '

Set-TextValue 12 4 '7.08'
Set-TextValue 12 5 '  -2.12%  '

Set-TextValue 13 5 '  -0.22%  '

Set-TextValue 14 4 '2.569.08'
Set-TextValue 14 5 '  -0.79%  '

Set-TextValue 15 4 '2.226.97'
Set-TextValue 15 5 '  -4.34%  '

Set-TextValue 16 4 '0.825'
Set-TextValue 16 5 '  -1.12%  '

Set-TextValue 17 4 '13.39'
Set-TextValue 17 5 '  -1.66%  '

Set-TextValue 18 4 '44.419.60'
Set-TextValue 18 5 '  +0.78%  '

Set-TextValue 19 4 '0.0₃0933'
Set-TextValue 19 5 '  -3.93%  '

Set-TextValue 20 2 'Uniswap'
Set-TextValue 20 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 20 4 '6.15'
Set-TextValue 20 5 '  -4.05%  '

Set-TextValue 21 2 'InternetComputer(DFINITY)'
Set-TextValue 21 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 21 4 '11.60'
Set-TextValue 21 5 '  -4.60%  '

Set-TextValue 22 4 '64.70'
Set-TextValue 22 5 '  -1.22%  '

Set-TextValue 23 4 '236.07'
Set-TextValue 23 5 '  -0.44%  '

Set-TextValue 24 4 '2.90'
Set-TextValue 24 5 '  -1.81%  '

Set-TextValue 25 4 '1.95'
Set-TextValue 25 5 '  -2.53%  '

Set-TextValue 26 5 '  -0.01%  '

Set-TextValue 27 4 '2.30'
Set-TextValue 27 5 '  +4.10%  '

Set-TextValue 28 4 '9.64'
Set-TextValue 28 5 '  -3.55%  '

Set-TextValue 29 4 '36.97'
Set-TextValue 29 5 '  -1.91%  '

Set-TextValue 30 2 'EthereumClassic'
Set-TextValue 30 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 30 4 '19.71'
Set-TextValue 30 5 '  -2.37%  '

Set-TextValue 31 2 'Filecoin'
Set-TextValue 31 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 31 4 '5.78'
Set-TextValue 31 5 '  -3.67%  '

Set-TextValue 32 4 '149.32'
Set-TextValue 32 5 '  -2.37%  '

Set-TextValue 33 2 'WEMIXToken'
Set-TextValue 33 3 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 33 4 '2.62'
Set-TextValue 33 5 '  +0.58%  '

Set-TextValue 34 2 'Hedera'
Set-TextValue 34 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 34 4 '0.0777'
Set-TextValue 34 5 '  -3.16%  '

Set-TextValue 35 5 '  -6.67%  '

Set-TextValue 36 5 '  -0.67%  '

Set-TextValue 37 4 '0.117'
Set-TextValue 37 5 '  -2.36%  '

Set-TextValue 38 4 '1.84'
Set-TextValue 38 5 '  +4.65%  '

Set-TextValue 39 4 '14.73'
Set-TextValue 39 5 '  +0.58%  '

Set-TextValue 40 4 '3.31'
Set-TextValue 40 5 '  -4.87%  '

Set-TextValue 41 4 '3.72'
Set-TextValue 41 5 '  -3.76%  '

Set-TextValue 42 4 '0.0295'
Set-TextValue 42 5 '  -0.76%  '

Set-TextValue 44 4 '1.816.35'
Set-TextValue 44 5 '  +4.09%  '

Set-TextValue 45 4 '1.72'
Set-TextValue 45 5 '  +8.42%  '

Set-TextValue 46 4 '78.55'
Set-TextValue 46 5 '  -5.57%  '

Set-TextValue 47 4 '0.185'
Set-TextValue 47 5 '  -3.37%  '

Set-TextValue 48 4 '97.33'
Set-TextValue 48 5 '  -2.97%  '

Set-TextValue 49 4 '4.80'
Set-TextValue 49 5 '  -2.95%  '

Set-TextValue 50 2 'FraxShare'
Set-TextValue 50 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 50 4 '7.92'
Set-TextValue 50 5 '  -2.94%  '

Set-TextValue 51 2 'ordi'
Set-TextValue 51 3 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue 51 4 '67.45'
Set-TextValue 51 5 '  -1.06%  '
